# Sync attendance_reports: reorder the "Recorded By" (column G) lists so
# that the comma-separated list of recorder identities is reversed for
# each row, except where it is exactly "admin@admin.com, System" which is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ",\s*"
    $n = $parts.Length

    if ($n -le 1) { continue }

    $skip = $false
    if ($n -eq 2 -and $parts[0] -eq "admin@admin.com" -and $parts[1] -eq "System") {
        $skip = $true
    }

    if ($skip) { continue }

    $reversedParts = @()
    for ($i = 0; $i -lt $n; $i++) {
        $idx = $n - 1 - $i
        $reversedParts += $parts[$idx]
    }

    $newVal = [string]::Join(", ", $reversedParts)
    $cell.Value = $newVal
}
